$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J2").Value = "Boolean"
$ws.Range("J3").Value = $false
$ws.Range("J3").Select() | Out-Null
